# The upstream analysis workflow was fixed so the sensitivity sweep now
# starts its cutoff window 4 steps later. For every data sheet, drop the
# first 4 data rows (old Cutoff = 0..3) and re-number the remaining rows'
# Cutoff column (A) back to a zero-based sequence, shifting the rest of
# the table up.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Remove the first four data rows (rows 2-5); everything below shifts up.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Re-number the Cutoff column (A) starting again from 0 for the rows
    # that remain (now in rows 2-16).
    for ($i = 0; $i -lt 15; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
